$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.903.10'
$ws.Range("E2").Value = '  +0.87%  '

$ws.Range("D3").Value = '2.288.58'
$ws.Range("E3").Value = '  -0.63%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.39'
$ws.Range("E5").Value = '  -0.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.63'
$ws.Range("E6").Value = '  +0.66%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -1.24%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.603'
$ws.Range("E9").Value = '  -1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.54'
$ws.Range("E10").Value = '  -1.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0902'
$ws.Range("E11").Value = '  -0.89%  '

$ws.Range("E12").Value = '  +0.59%  '

$ws.Range("E13").Value = '  +2.27%  '

$ws.Range("E14").Value = '  +3.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.23'
$ws.Range("E15").Value = '  -0.77%  '

$ws.Range("D16").Value = '2.635.91'
$ws.Range("E16").Value = '  -0.68%  '

$ws.Range("D17").Value = '2.288.85'
$ws.Range("E17").Value = '  -0.63%  '

$ws.Range("D18").Value = '42.598.90'
$ws.Range("E18").Value = '  +0.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.42'
$ws.Range("E19").Value = '  -0.69%  '

$ws.Range("E20").Value = '  -0.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.58'
$ws.Range("E21").Value = '  +22.34%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.91'
$ws.Range("E22").Value = '  +0.75%  '

$ws.Range("E23").Value = '  +0.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '262.48'
$ws.Range("E24").Value = '  -4.90%  '

$ws.Range("E25").Value = '  -3.34%  '

$ws.Range("E26").Value = '  +0.51%  '

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.16'
$ws.Range("E28").Value = '  +21.69%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.29'
$ws.Range("E29").Value = '  -2.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.28'
$ws.Range("E30").Value = '  -2.23%  '

$ws.Range("E31").Value = '  +4.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '166.59'
$ws.Range("E32").Value = '  +0.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0874'
$ws.Range("E33").Value = '  +0.25%  '

$ws.Range("E34").Value = '  -3.65%  '

$ws.Range("E35").Value = '  -0.68%  '

$ws.Range("E36").Value = '  -2.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.55'
$ws.Range("E37").Value = '  -1.17%  '

$ws.Range("E38").Value = '  -5.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.83'
$ws.Range("E39").Value = '  +2.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.67'
$ws.Range("E40").Value = '  -3.15%  '

$ws.Range("E41").Value = '  +5.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.232'
$ws.Range("E42").Value = '  +1.93%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.72'
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("E44").Value = '  -0.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.39'
$ws.Range("E45").Value = '  -2.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.21'
$ws.Range("E46").Value = '  +1.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '113.55'
$ws.Range("E47").Value = '  +0.28%  '

$ws.Range("D48").Value = '1.727.02'
$ws.Range("E48").Value = '  +8.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '79.27'
$ws.Range("E49").Value = '  -3.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.77'
$ws.Range("E50").Value = '  -1.67%  '

$ws.Range("E51").Value = '  +0.41%  '
